$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

$replacements = @(
    @(1,1,"34÷6="),
    @(1,2,"42÷6="),
    @(1,3,"44÷9="),
    @(1,4,"44÷9="),
    @(1,5,"45÷9="),
    @(5,1,"71÷4="),
    @(5,2,"18÷6="),
    @(5,3,"62÷2="),
    @(5,4,"43÷8="),
    @(5,5,"50÷4="),
    @(9,1,"52÷7="),
    @(9,2,"42÷9="),
    @(9,3,"14÷9="),
    @(9,4,"67÷8="),
    @(9,5,"70÷6="),
    @(13,1,"88÷7="),
    @(13,2,"73÷8="),
    @(13,3,"64÷7="),
    @(13,4,"44÷8="),
    @(13,5,"31÷2="),
    @(17,1,"31÷3="),
    @(17,2,"49÷2="),
    @(17,3,"47÷3="),
    @(17,4,"45÷6="),
    @(17,5,"76÷2=")
)

foreach ($r in $replacements) {
    $row = $r[0]
    $col = $r[1]
    $newText = $r[2]
    $cell = $tbl.Cell($row, $col)
    $rng = $cell.Range
    $rng.End = $rng.End - 1
    $rng.Text = $newText
}
